# excel persona names update
#
# The analyst added a "cluster_name" column (column C) giving each persona
# cluster a human first name, and tidied up two "highly educated" labels to
# just "educated" in the description text. Selection/scroll position was
# also reset after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New persona names in column C (cluster_name), one per cluster row.
$ws.Range("C4").Value = "Susana"
$ws.Range("C5").Value = "Richard"
$ws.Range("C6").Value = "Jennifer"
$ws.Range("C7").Value = "Michael"
$ws.Range("C8").Value = "Karen"

# "highly educated" -> "educated" wording tweak in the affected descriptions.
$ws.Range("D5").Value = " educated , older couples / small families"
$ws.Range("AB5").Value = "family oriented,  educated, medium income, wine, online, deals sometimes, teenager, 40s"

# Leave the view scrolled back to the left with the selection on D15.
$ws.Range("D15").Select()
